$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Ex 1-A")
$r = $ws1.Range("A17")
$r.Value = "Question"
$r.Font.Bold = $true

$r2 = $ws1.Range("A18")
$r2.Value = "For which value of S&P do we get a output of the same dimension as the input?"

$r3 = $ws1.Range("A19")
$r3.Value = "S1P1"
$r3.Font.Bold = $true
